# Y4_B2526_General_Surgery scanner log - edited session update.
# Source diff:
#   1. Worksheet tab renamed General_Surgery -> Session.
#   2. A stray duplicate/out-of-order log line (Student ID 201834 @ 11:35:22,
#      which had been inserted ahead of where it belonged) is removed. Every
#      row below it shifts up by one, so the sheet's used range shrinks from
#      A1:F114 to A1:F113.
#
# Deleting the entire worksheet row reproduces exactly that: Excel shifts
# every lower row up by one and shrinks the sheet dimension automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra log entry that was recorded as row 41
# (Student ID 201834, logged 11:35:22) - everything beneath shifts up.
$ws.Rows.Item(41).Delete()

# Rename the worksheet tab to match the new session name.
$ws.Name = "Session"
